# agent integration and initial UI created
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clean up the question text in A11 (drop the internal testing note)
$ws.Range("A11").Value = "Are birds green?"

# Update the active selection on Sheet1 to A2:A10 (anchor A2)
$ws.Activate()
$ws.Range("A2:A10").Select()
